$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: car model names now wrap onto multiple lines (space -> line break)
$ws.Range("B1").Value = "Kia" + [char]10 + "Rio"
$ws.Range("C1").Value = "Volkswagen" + [char]10 + "Golf"
$ws.Range("D1").Value = "Toyota" + [char]10 + "Corolla"
$ws.Range("E1").Value = "Skoda" + [char]10 + "Octavia"
$ws.Range("F1").Value = "BMW" + [char]10 + "3" + [char]10 + "Series"
$ws.Range("G1").Value = "Hyundai" + [char]10 + "Solaris"
$ws.Range("H1").Value = "Вектор" + [char]10 + "приоритетов"

# Column A (row labels): same multi-line treatment
$ws.Range("A2").Value = "Kia" + [char]10 + "Rio"
$ws.Range("A3").Value = "Volkswagen" + [char]10 + "Golf"
$ws.Range("A4").Value = "Toyota" + [char]10 + "Corolla"
$ws.Range("A5").Value = "Skoda" + [char]10 + "Octavia"
$ws.Range("A6").Value = "BMW" + [char]10 + "3" + [char]10 + "Series"
$ws.Range("A7").Value = "Hyundai" + [char]10 + "Solaris"

# Recomputed priority-vector values (kept as text, same as the rest of column H)
$ws.Range("H4").Value = "'0.228"
$ws.Range("H6").Value = "'0.362"

# Columns got narrower now that the headers/labels wrap across multiple lines
$ws.Columns.Item(1).ColumnWidth = 14.4 - 5/6
$ws.Columns.Item(2).ColumnWidth = 6 - 5/6
$ws.Columns.Item(3).ColumnWidth = 14.4 - 5/6
$ws.Columns.Item(4).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(5).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(6).ColumnWidth = 9.6 - 5/6
$ws.Columns.Item(7).ColumnWidth = 10.8 - 5/6
$ws.Columns.Item(8).ColumnWidth = 15.6 - 5/6
